# Apply the data edits from the diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2: 4 -> 4.5
$ws.Range("B2").Value = 4.5

# C4: 1.45 -> 1.4
$ws.Range("C4").Value = 1.4

# C5: 20 -> 25
$ws.Range("C5").Value = 25

# Move the active selection from C5 to C4 (matches saved cursor position).
$ws.Range("C4").Select()
